$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column I (same rows) onto the new column J
# before writing values, so the new cells pick up style index 1.
$ws.Range("I1:I4").Copy()
$ws.Range("J1:J4").PasteSpecial(-4122)

# Add new "attraction" column (J) header + data
$ws.Range("J1").Value = "attraction"
$ws.Range("J2").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("J4").Value = 5

# Update the active selection to match the target state
$ws.Range("J5").Select()
